# Commit: "update to 2024 data"
#
# Adds a 4th sow-fertilization event (date 2024-04-23, npk 2-2-2,
# "Espoma- Organic Grow!", dose "half") to the crop log:
#   - crops 2-8   (already had an explicit-but-empty 3rd fert event in W:Z):
#       only need the previously style-only-blank 4th event columns (AA:AD)
#       turned into explicit empty text, matching W:Z's empty text.
#   - crops 9-20  (already had a real 3rd fert event in W:Z):
#       the 4th event columns (AA:AD) get the real new fert-event values.
#   - crops 21-33 (still had an empty 3rd fert event in W:Z):
#       the 3rd event columns (W:Z) get the new fert-event values (becoming
#       the first fert event recorded for these crops this round), and the
#       4th event columns (AA:AD) become explicit empty text (mirrors 2-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    # Force text storage so date-shaped strings ("2024-04-23") are not
    # auto-parsed into Excel date serials, then drop the resulting style
    # delta so the cell keeps the workbook's default (styleless) look.
    if ($val -match '^\d{4}-\d{2}-\d{2}$') {
        $ws.Range($addr).NumberFormat = "@"
    }
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

$fertDate = "2024-04-23"
$fertNpk = "2-2-2"
$fertName = "Espoma- Organic Grow!"
$fertDose = "half"

# Rows 2-8: garlic, kale, spinach, lettuce x4
for ($r = 2; $r -le 8; $r++) {
    Set-TextValue $ws "AA$r" ""
    Set-TextValue $ws "AB$r" ""
    Set-TextValue $ws "AC$r" ""
    Set-TextValue $ws "AD$r" ""
}

# Rows 9-20: eggplant, pepper varieties
for ($r = 9; $r -le 20; $r++) {
    Set-TextValue $ws "AA$r" $fertDate
    Set-TextValue $ws "AB$r" $fertNpk
    Set-TextValue $ws "AC$r" $fertName
    Set-TextValue $ws "AD$r" $fertDose
}

# Rows 21-33: tomato varieties, broccoli, basil, catnip, cilantro x2, coleus,
# marigold, coreopsis, salvia
for ($r = 21; $r -le 33; $r++) {
    Set-TextValue $ws "W$r" $fertDate
    Set-TextValue $ws "X$r" $fertNpk
    Set-TextValue $ws "Y$r" $fertName
    Set-TextValue $ws "Z$r" $fertDose
    Set-TextValue $ws "AA$r" ""
    Set-TextValue $ws "AB$r" ""
    Set-TextValue $ws "AC$r" ""
    Set-TextValue $ws "AD$r" ""
}
